# The workbook starts with a single worksheet ("hclust") holding the
# cluster data. Bring it in line with the re-saved (STATA 18) layout:
#   - rename "hclust" -> "Sheet1" (keeps all of its data)
#   - append two new, empty worksheets: "Sheet2" and "Sheet3"
#   - select the whole data sheet (A1:XFD1048576), then make "Sheet2"
#     the active/visible tab, as left behind in the saved file.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sheet1"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws3 = $wb.Worksheets.Add($null, $ws2)

$ws1.Cells.Select() | Out-Null
$ws2.Activate() | Out-Null
